# Apply the "synthetic_constrained_1" workbook edit:
# - clear the now-"undefined" N30 cell (so AVERAGE(N2:N31) below ignores it)
# - add summary rows: a plain average-of-J row, two labeled averages
#   (SW and SC ratios) and two labeled "worst" (min/max) literals
# - format the new summary cells with bold fonts, matching the source

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The optimum was undefined for dataset row 30 (same situation as row 25),
# so N30 is cleared to blank instead of holding a literal 0.
$ws.Range("N30").ClearContents()

# Row 32: plain average of the |S*|/n column (J)
$ws.Range("J32").Formula = "=AVERAGE(J2:J31)"
$ws.Range("J32").Font.Bold = $true
$ws.Range("J32").Font.Size = 12
$ws.Range("J32").VerticalAlignment = -4108
$ws.Rows.Item(32).RowHeight = 15.6

# Row 34-35: labeled averages of the SW/SC ratio columns (N and Z)
$ws.Range("A34").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B34").Formula = "=AVERAGE(N2:N31)"
$ws.Range("B34").Font.Bold = $true
$ws.Range("B34").Font.Size = 12
$ws.Range("B34").VerticalAlignment = -4108
$ws.Rows.Item(34).RowHeight = 15.6

$ws.Range("A35").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B35").Formula = "=AVERAGE(Z2:Z31)"
$ws.Range("B35").Font.Bold = $true
$ws.Range("B35").Font.Size = 12
$ws.Range("B35").VerticalAlignment = -4108
$ws.Rows.Item(35).RowHeight = 15.6

# Row 36-37: labeled worst-case (min/max) values, entered as literals
$ws.Range("A36").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B36").Value = 0.62075470867252491
$ws.Range("B36").Font.Bold = $true
$ws.Range("B36").Font.Size = 11

$ws.Range("A37").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B37").Value = 1.1605844883682661
$ws.Range("B37").Font.Bold = $true
$ws.Range("B37").Font.Size = 11

# Selection / print setup, matching the saved view state
$ws.Range("A34:B37").Select() | Out-Null
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
